# Re-saved/re-downloaded weather data on 2025-05-11 13:16 — refresh the
# "Nedladdat" timestamp, the DATUM date, and the Hour/Temperature/Rainorsnow
# readings for every data row (2-27) to the new snapshot's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$downloaded = "2025-05-11 13:16"

# row -> @(DATUM, Hour, Temperature, Rainorsnow)
$rows = @{
    2  = @("2025-05-11", 13, 12.4, $false)
    3  = @("2025-05-11", 14, 13.3, $false)
    4  = @("2025-05-11", 15, 12.5, $false)
    5  = @("2025-05-11", 16, 12.4, $false)
    6  = @("2025-05-11", 17, 12,   $false)
    7  = @("2025-05-11", 18, 11.2, $false)
    8  = @("2025-05-11", 19, 9.6,  $false)
    9  = @("2025-05-11", 20, 8,    $false)
    10 = @("2025-05-11", 21, 6.8,  $false)
    11 = @("2025-05-11", 22, 6.1,  $false)
    12 = @("2025-05-11", 23, 5.9,  $false)
    13 = @("2025-05-12", 0,  5.4,  $false)
    14 = @("2025-05-12", 1,  4.7,  $false)
    15 = @("2025-05-12", 2,  4.6,  $false)
    16 = @("2025-05-12", 3,  4.3,  $false)
    17 = @("2025-05-12", 4,  6.1,  $false)
    18 = @("2025-05-12", 5,  9,    $false)
    19 = @("2025-05-12", 6,  10.3, $false)
    20 = @("2025-05-12", 7,  11.1, $false)
    21 = @("2025-05-12", 8,  11.7, $false)
    22 = @("2025-05-12", 9,  11.9, $false)
    23 = @("2025-05-12", 10, 12.7, $false)
    24 = @("2025-05-12", 11, 14.2, $false)
    25 = @("2025-05-12", 12, 14.9, $false)
    26 = @("2025-05-12", 13, 14.9, $false)
    27 = @("2025-05-12", 14, 15.1, $false)
}

foreach ($r in 2..27) {
    $datum = $rows[$r][0]
    $hour = $rows[$r][1]
    $temp = $rows[$r][2]
    $rainorsnow = $rows[$r][3]

    $ws.Range("A$r").Value = $downloaded

    # Force DATUM to stay a plain text cell (leading apostrophe prevents
    # Excel from auto-converting the yyyy-mm-dd string into a date serial),
    # then drop the resulting quote-prefix style so the cell stays unstyled.
    $ws.Range("D$r").Value = "'" + $datum
    $ws.Range("D$r").ClearFormats()

    $ws.Range("E$r").Value = $hour
    $ws.Range("F$r").Value = $temp
    $ws.Range("G$r").Value = $rainorsnow
}
